$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(3)
$sh.TextFrame.TextRange.Font.Bold = 1
Write-Output "done"
